$d = $word.ActiveDocument

$replacements = @(
    @("18×68=1224", "49×77=3773"),
    @("68×98=6664", "96×54=5184"),
    @("49×69=3381", "89×17=1513"),
    @("51×99=5049", "77×68=5236"),
    @("34×33=1122", "91×48=4368"),
    @("17×81=1377", "77×21=1617"),
    @("65×39=2535", "46×61=2806"),
    @("68×31=2108", "37×81=2997"),
    @("19×35=665",  "94×37=3478"),
    @("18×60=1080", "30×20=600"),
    @("18×25=450",  "32×19=608"),
    @("26×27=702",  "46×88=4048"),
    @("55×41=2255", "32×44=1408"),
    @("35×58=2030", "88×74=6512"),
    @("52×23=1196", "77×83=6391"),
    @("42×89=3738", "55×61=3355"),
    @("40×52=2080", "43×14=602"),
    @("89×81=7209", "90×74=6660"),
    @("65×34=2210", "84×68=5712"),
    @("68×18=1224", "36×63=2268"),
    @("48×32=1536", "19×92=1748"),
    @("38×65=2470", "65×85=5525"),
    @("18×66=1188", "89×33=2937"),
    @("18×47=846",  "20×76=1520"),
    @("43×15=645",  "45×94=4230")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
